$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-26 17:48:43'
$ws.Range('E3').Value = '2026-02-26 17:48:45'
$ws.Range('K3').Value = '17.2 MJ/m2'
$ws.Range('E4').Value = '2026-02-26 17:48:48'
$ws.Range('O4').Value = '11.1 °C'
$ws.Range('E5').Value = '2026-02-26 17:48:51'
$ws.Range('H5').Value = '''44%'
$ws.Range('K5').Value = '16.7 MJ/m2'
$ws.Range('O5').Value = '5.3 °C'
$ws.Range('E6').Value = '2026-02-26 17:48:54'
$ws.Range('E7').Value = '2026-02-26 17:48:57'
$ws.Range('K7').Value = '15.5 MJ/m2'
$ws.Range('E8').Value = '2026-02-26 17:48:59'
$ws.Range('K8').Value = '15.2 MJ/m2'
$ws.Range('O8').Value = '11.5 °C'
$ws.Range('E9').Value = '2026-02-26 17:49:02'
$ws.Range('E10').Value = '2026-02-26 17:49:05'
$ws.Range('K10').Value = '14.6 MJ/m2'
$ws.Range('E11').Value = '2026-02-26 17:49:07'
$ws.Range('O11').Value = '8.8 °C'
$ws.Range('E12').Value = '2026-02-26 17:49:10'
$ws.Range('O12').Value = '11.7 °C'
$ws.Range('E13').Value = '2026-02-26 17:49:12'
$ws.Range('H13').Value = '''64%'
$ws.Range('J13').Value = '1028.4 hPa'
$ws.Range('O13').Value = '6.9 °C'
$ws.Range('E14').Value = '2026-02-26 17:49:15'
$ws.Range('K14').Value = '15.1 MJ/m2'
$ws.Range('E15').Value = '2026-02-26 17:49:18'
$ws.Range('E16').Value = '2026-02-26 17:49:21'
$ws.Range('K16').Value = '15.3 MJ/m2'
$ws.Range('E17').Value = '2026-02-26 17:49:23'
$ws.Range('K17').Value = '17.9 MJ/m2'
$ws.Range('E18').Value = '2026-02-26 17:49:26'
$ws.Range('O18').Value = '12.1 °C'
$ws.Range('E19').Value = '2026-02-26 17:49:29'
$ws.Range('K19').Value = '16.0 MJ/m2'
$ws.Range('E20').Value = '2026-02-26 17:49:32'
$ws.Range('H20').Value = '''46%'
$ws.Range('K20').Value = '17.2 MJ/m2'
$ws.Range('O20').Value = '3.0 °C'
$ws.Range('E21').Value = '2026-02-26 17:49:34'
$ws.Range('H21').Value = '''63%'
$ws.Range('J21').Value = '1027.1 hPa'
$ws.Range('O21').Value = '9.8 °C'
$ws.Range('E22').Value = '2026-02-26 17:49:37'
$ws.Range('K22').Value = '17.6 MJ/m2'
$ws.Range('E23').Value = '2026-02-26 17:49:40'
$ws.Range('K23').Value = '17.0 MJ/m2'
$ws.Range('O23').Value = '3.5 °C'
$ws.Range('E24').Value = '2026-02-26 17:49:43'
$ws.Range('H24').Value = '''73%'
$ws.Range('K24').Value = '16.1 MJ/m2'
$ws.Range('O24').Value = '10.6 °C'
$ws.Range('E25').Value = '2026-02-26 17:49:46'
$ws.Range('E26').Value = '2026-02-26 17:49:48'
$ws.Range('H26').Value = '''37%'
$ws.Range('E27').Value = '2026-02-26 17:49:51'
$ws.Range('H27').Value = '''40%'
$ws.Range('O27').Value = '5.3 °C'
$ws.Range('E28').Value = '2026-02-26 17:49:53'
$ws.Range('H28').Value = '''76%'
$ws.Range('K28').Value = '14.1 MJ/m2'
$ws.Range('O28').Value = '11.1 °C'
$ws.Range('E29').Value = '2026-02-26 17:49:56'
$ws.Range('E30').Value = '2026-02-26 17:49:59'
$ws.Range('K30').Value = '15.1 MJ/m2'
$ws.Range('E31').Value = '2026-02-26 17:50:02'
$ws.Range('K31').Value = '14.5 MJ/m2'
$ws.Range('E32').Value = '2026-02-26 17:50:05'
$ws.Range('O32').Value = '8.5 °C'
$ws.Range('E33').Value = '2026-02-26 17:50:07'
$ws.Range('J33').Value = '1026.8 hPa'
$ws.Range('O33').Value = '8.5 °C'
$ws.Range('E34').Value = '2026-02-26 17:50:09'
$ws.Range('K34').Value = '15.2 MJ/m2'
$ws.Range('O34').Value = '5.3 °C'
$ws.Range('E35').Value = '2026-02-26 17:50:12'
$ws.Range('K35').Value = '16.8 MJ/m2'
$ws.Range('E36').Value = '2026-02-26 17:50:15'
$ws.Range('E37').Value = '2026-02-26 17:50:18'
$ws.Range('J37').Value = '1028.0 hPa'
$ws.Range('O37').Value = '8.1 °C'
$ws.Range('E38').Value = '2026-02-26 17:50:20'
$ws.Range('K38').Value = '14.9 MJ/m2'
$ws.Range('O38').Value = '11.5 °C'
$ws.Range('E39').Value = '2026-02-26 17:50:23'
$ws.Range('O39').Value = '3.1 °C'
$ws.Range('E40').Value = '2026-02-26 17:50:26'
$ws.Range('H40').Value = '''66%'
$ws.Range('J40').Value = '1027.4 hPa'
$ws.Range('O40').Value = '9.7 °C'
$ws.Range('E41').Value = '2026-02-26 17:50:28'
$ws.Range('K41').Value = '16.0 MJ/m2'
$ws.Range('O41').Value = '11.4 °C'
$ws.Range('E42').Value = '2026-02-26 17:50:31'
$ws.Range('O42').Value = '11.6 °C'
$ws.Range('E43').Value = '2026-02-26 17:50:34'
$ws.Range('K43').Value = '15.8 MJ/m2'
$ws.Range('O43').Value = '9.3 °C'
$ws.Range('E44').Value = '2026-02-26 17:50:36'
$ws.Range('E45').Value = '2026-02-26 17:50:39'
$ws.Range('J45').Value = '1025.9 hPa'
$ws.Range('K45').Value = '14.9 MJ/m2'
$ws.Range('O45').Value = '10.8 °C'
$ws.Range('E46').Value = '2026-02-26 17:50:42'
$ws.Range('K46').Value = '12.8 MJ/m2'
$ws.Range('O46').Value = '11.1 °C'
